# Generate Report for Handback
# Row 3 (eb5fca58-5717-4c65-9b3d-2ba88abb2acc.md) has now been handed back
# (in sync with en-US), so its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet that lists it, and the
# per-locale "Latest Handback DateTime" is stamped with the handback time.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-03-11 03:11:43"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-03-11 03:12:04"
